# Adm-Calcr.xlsx — refresh LR-pairs table with the new TPM run.
#
# The sending/target cluster set changes ("Inflammatory-Mac" is replaced by
# "Resolving-Mac", elsewhere in the cluster list), every numeric column is
# recomputed, and the table grows from 5 data rows to 10 data rows (2 target
# clusters x 5 sending clusters), so the sheet's used range grows from
# A1:T6 to A1:T11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol
#          D Target cluster  | E..T numeric NATMI edge-scoring stats
$rows = @(
  ,@("ECs",           "Adm", "Calcr", "ECs",    3, 1,                   15.780133,            47.340399,  0.3664158323494666,   0.3664158323494666,   1, 0.3333333333333333, 0.09584066666666667, 0.287522,            0.3176649051884244, 0.3176649051884244, 1.512378466808667,   13.611406201278,    0.1163974506428309,   0.1163974506428309)
  ,@("ECs",           "Adm", "Calcr", "MuSCs",  3, 1,                   15.780133,            47.340399,  0.3664158323494666,   0.3664158323494666,   2, 0.6666666666666666, 0.205863,            0.6175889999999999, 0.6823350948115756, 0.6823350948115755, 3.248545519779,      29.236909678011,    0.2500183817066357,   0.2500183817066356)
  ,@("FAPs",          "Adm", "Calcr", "ECs",    3, 1,                   26.269711,            78.809133,  0.6099845940236968,   0.6099845940236966,   1, 0.3333333333333333, 0.09584066666666667, 0.287522,            0.3176649051884244, 0.3176649051884244, 2.517706615380667,   22.659359538426,    0.1937706982269372,   0.1937706982269371)
  ,@("FAPs",          "Adm", "Calcr", "MuSCs",  3, 1,                   26.269711,            78.809133,  0.6099845940236968,   0.6099845940236966,   2, 0.6666666666666666, 0.205863,            0.6175889999999999, 0.6823350948115756, 0.6823350948115755, 5.407961515593,      48.671653640337,    0.4162138957967596,   0.4162138957967594)
  ,@("MuSCs",         "Adm", "Calcr", "ECs",    3, 1,                   0.973204,             2.919612,   0.02259786743912933,  0.02259786743912933,  1, 0.3333333333333333, 0.09584066666666667, 0.287522,            0.3176649051884244, 0.3176649051884244, 0.09327252016266667, 0.839452681464,     0.007178549417511604, 0.007178549417511601)
  ,@("MuSCs",         "Adm", "Calcr", "MuSCs",  3, 1,                   0.973204,             2.919612,   0.02259786743912933,  0.02259786743912933,  2, 0.6666666666666666, 0.205863,            0.6175889999999999, 0.6823350948115756, 0.6823350948115755, 0.200346695052,      1.803120255468,     0.01541931802161773,  0.01541931802161773)
  ,@("Neutrophils",   "Adm", "Calcr", "ECs",    1, 0.3333333333333333,  0.01491966666666667, 0.044759,    0.0003464357417040312,0.0003464357417040312,1, 0.3333333333333333, 0.09584066666666667, 0.287522,            0.3176649051884244, 0.3176649051884244, 0.001429910799777778,0.012869197198,     0.0001100504770422926,0.0001100504770422925)
  ,@("Neutrophils",   "Adm", "Calcr", "MuSCs",  1, 0.3333333333333333,  0.01491966666666667, 0.044759,    0.0003464357417040312,0.0003464357417040312,2, 0.6666666666666666, 0.205863,            0.6175889999999999, 0.6823350948115756, 0.6823350948115755, 0.003071407339,      0.027642666051,     0.0002363852646617386,0.0002363852646617386)
  ,@("Resolving-Mac", "Adm", "Calcr", "ECs",    1, 0.3333333333333333,  0.02822,              0.08466,    0.0006552704460033352,0.0006552704460033352,1, 0.3333333333333333, 0.09584066666666667, 0.287522,            0.3176649051884244, 0.3176649051884244, 0.002704623613333334,0.02434161252,      0.0002081564241024261,0.000208156424102426)
  ,@("Resolving-Mac", "Adm", "Calcr", "MuSCs",  1, 0.3333333333333333,  0.02822,              0.08466,    0.0006552704460033352,0.0006552704460033352,2, 0.6666666666666666, 0.205863,            0.6175889999999999, 0.6823350948115756, 0.6823350948115755, 0.005809453859999999,0.05228508473999999,0.0004471140219009092,0.0004471140219009091)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2   # data starts on row 2 (row 1 is the header)
  $row = $rows[$i]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}
